# Logboek update: add presentatie / infoactivity entries for 14-12-2015 (maandag)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the hour values first so the SUM(D:D) dependents (F2/G2) pick them up
# cleanly on recalc.
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 0.5

# Copy formatting of the last existing data row (21) down onto the two new rows (22, 23)
$ws.Range("A21:D21").Copy() | Out-Null
$ws.Range("A22:D22").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:D23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 22: online meeting, betere taakverdeling, en wekelijkse resultaten
$ws.Range("A22").Value = 42352
$ws.Range("B22").Value = "maandag"
$ws.Range("C22").Value = "online meeting, betere taakverdeling, en wekelijkse resultaten"

# Row 23: Kleine powerpoint, infoactivity afmaken
$ws.Range("A23").Value = 42352
$ws.Range("B23").Value = "maandag"
$ws.Range("C23").Value = "Kleine powerpoint, infoactivity afmaken"

# Match the author's final selection in the saved view
$ws.Range("A24").Select() | Out-Null
